# lines_states.xlsx: add two new "line" entries (line7, line8) after line6,
# and two new "extr" entries (extr7, extr8) at the end, renumbering the
# index column (A) and updating a couple of C/D/E values along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data grid for rows 2..17 (A index, B name, C, D, E in_service)
$data = @(
    @(0,  "line1", 7,  9,  $true),
    @(1,  "line2", 9,  8,  $false),
    @(2,  "line3", 8,  10, $true),
    @(3,  "line4", 8,  11, $true),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $true),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $false),
    @(11, "extr4", 7,  8,  $false),
    @(12, "extr5", 9,  11, $true),
    @(13, "extr6", 7,  11, $false),
    @(14, "extr7", 5,  7,  $true),
    @(15, "extr8", 8,  5,  $false)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Two new rows (16, 17) were added at the bottom; give column A there the
# same bold/border/centered style ("s=1") the rest of the index column uses,
# by copying the format from the row above.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
